{"js": "// Replace the 25 division-problem strings (in document order) inside the\n// single worksheet table. The table has 20 rows x 5 columns; only every\n// 5th row (0, 4, 8, 12, 16) actually holds text, the others are blank\n// spacer rows. We read the table's 2D `values`, and replace each\n// non-empty cell's text with the next value from the ordered replacement\n// list below (same left-to-right, top-to-bottom order as the source\n// document), then write the grid back.\n\nconst replacements = [\n  \"10\u00f75=2, 0\",\n  \"40\u00f79=4, 4\",\n  \"53\u00f74=13, 1\",\n  \"91\u00f77=13, 0\",\n  \"84\u00f76=14, 0\",\n  \"23\u00f79=2, 5\",\n  \"73\u00f74=18, 1\",\n  \"92\u00f77=13, 1\",\n  \"27\u00f79=3, 0\",\n  \"54\u00f75=10, 4\",\n  \"91\u00f73=30, 1\",\n  \"99\u00f78=12, 3\",\n  \"46\u00f72=23, 0\",\n  \"92\u00f79=10, 2\",\n  \"90\u00f73=30, 0\",\n  \"13\u00f79=1, 4\",\n  \"81\u00f74=20, 1\",\n  \"11\u00f74=2, 3\",\n  \"65\u00f73=21, 2\",\n  \"45\u00f76=7, 3\",\n  \"52\u00f75=10, 2\",\n  \"54\u00f78=6, 6\",\n  \"56\u00f75=11, 1\",\n  \"51\u00f79=5, 6\",\n  \"21\u00f75=4, 1\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst grid = table.values;\nlet idx = 0;\nconst newGrid = grid.map((row) =>\n  row.map((cell) => {\n    if (cell !== \"\" && idx < replacements.length) {\n      return replacements[idx++];\n    }\n    return cell;\n  })\n);\n\ntable.values = newGrid;\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem strings (in document order) inside the\n# single worksheet table. The table is 20 rows x 5 columns; only every\n# 5th row (1, 5, 9, 13, 17 in 1-based indexing) actually holds text, the\n# others are blank spacer rows. Walk every cell in row-major order and\n# overwrite the non-empty ones with the next value from the ordered\n# replacement list below (same left-to-right, top-to-bottom order as the\n# source document).\n\n$replacements = @(\n    \"10\u00f75=2, 0\",\n    \"40\u00f79=4, 4\",\n    \"53\u00f74=13, 1\",\n    \"91\u00f77=13, 0\",\n    \"84\u00f76=14, 0\",\n    \"23\u00f79=2, 5\",\n    \"73\u00f74=18, 1\",\n    \"92\u00f77=13, 1\",\n    \"27\u00f79=3, 0\",\n    \"54\u00f75=10, 4\",\n    \"91\u00f73=30, 1\",\n    \"99\u00f78=12, 3\",\n    \"46\u00f72=23, 0\",\n    \"92\u00f79=10, 2\",\n    \"90\u00f73=30, 0\",\n    \"13\u00f79=1, 4\",\n    \"81\u00f74=20, 1\",\n    \"11\u00f74=2, 3\",\n    \"65\u00f73=21, 2\",\n    \"45\u00f76=7, 3\",\n    \"52\u00f75=10, 2\",\n    \"54\u00f78=6, 6\",\n    \"56\u00f75=11, 1\",\n    \"51\u00f79=5, 6\",\n    \"21\u00f75=4, 1\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$idx = 0\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cellText = $cell.Range.Text\n        $cellText = $cellText.TrimEnd([char]7).TrimEnd([char]13)\n        if ($cellText -ne \"\" -and $idx -lt $replacements.Length) {\n            $cell.Range.Text = $replacements[$idx]\n            $idx = $idx + 1\n        }\n    }\n}\n"}
